$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - path gains a leading slash (rest unchanged)
$ws.Range("A2").Value = "/about-cancer/coping/feelings"

# Row 3 - now the Spanish "relax" page data (replaces the old Spanish "feelings" row)
$ws.Range("A3").Value = "/espanol/cancer/sobrellevar/sentimientos/relajarse"
$ws.Range("B3").Value = "Article"
$ws.Range("C3").Value = "Spanish"
$ws.Range("D3").Value = "Aprenda a relajarse | CGDP - Dev"
$ws.Range("E3").Value = "Aprenda a relajarse"

# Row 4 stays "Learning to Relax" (English) - unchanged from before

# Row 5 - new "Dana Farber/Harvard Cancer Center" row (replaces old Spanish "Aprenda a relajarse" row)
$ws.Range("A5").Value = "/about-cancer/coping/feelings/relaxation/dfharvard"
$ws.Range("B5").Value = "Cancer Center"
$ws.Range("C5").Value = "English"
$ws.Range("D5").Value = "Dana Farber/Harvard Cancer Center | CGDP - Dev"
$ws.Range("E5").Value = "Dana Farber/Harvard Cancer Center"
$ws.Range("E5").Font.Color = 2236962

# Row 6 - new "Jennifer K. Loukissas" biography row (replaces old "Duke Cancer Center" row)
$ws.Range("A6").Value = "/about-cancer/coping/feelings/relaxation/loukissas-jennifer"
$ws.Range("B6").Value = "Biography"
$ws.Range("C6").Value = "English"
$ws.Range("D6").Value = "Jennifer K. Loukissas, M.P.P. | CGDP - Dev"
$ws.Range("E6").Value = "Jennifer K. Loukissas, M.P.P."

# Column widths / best-fit (closest achievable values given this host's pixel rounding)
$ws.Columns("A").ColumnWidth = 36
$ws.Columns("B").ColumnWidth = 12.5
$ws.Columns("D").ColumnWidth = 35
$ws.Columns("E").ColumnWidth = 12.5

# Page setup
$ws.PageSetup.Orientation = 1

# Selection
[void]$ws.Range("H12").Select()
